$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new fixture rows above the existing data (rows 2 and 3);
# everything that was row 2 onward shifts down to row 4 onward.
$ws.Rows("2:3").Insert()
# Insert() copies formatting from the row above (the bold header) - strip it
# so the new rows match the plain (unstyled) formatting of the other data rows.
$ws.Rows("2:3").ClearFormats()

# --- New row 2: Gillingham vs Bromley ---
$ws.Range("A2").Value = 'Základní Poisson'
$ws.Range("B2").Value = 'e3'
$ws.Range("C2").Value = 'Gillingham'
$ws.Range("D2").Value = 'Bromley'
$ws.Range("E2").Value = 1480.630280430628
$ws.Range("F2").Value = 1532.144152754764
$ws.Range("G2").Value = "'6.02"
$ws.Range("H2").Value = "'4.10"
$ws.Range("I2").Value = "'1.72"
$ws.Range("J2").Value = "'2.52"
$ws.Range("K2").Value = "'1.66"
$ws.Range("L2").Value = 'Bivariantní Poisson'
$ws.Range("M2").Value = "'3.07"
$ws.Range("N2").Value = "'3.16"
$ws.Range("O2").Value = "'2.79"
$ws.Range("P2").Value = "'2.25"
$ws.Range("Q2").Value = "'1.80"
$ws.Range("R2").Value = 'Monte Carlo'
$ws.Range("S2").Value = "'2.80"
$ws.Range("T2").Value = "'3.74"
$ws.Range("U2").Value = "'2.67"
$ws.Range("V2").Value = "'2.23"
$ws.Range("W2").Value = "'1.81"
$ws.Range("X2").Value = 'xGBoost'
$ws.Range("Y2").Value = "'5.51"
$ws.Range("Z2").Value = "'1.31"
$ws.Range("AA2").Value = "'18.19"
$ws.Range("AB2").Value = "'6.07"
$ws.Range("AC2").Value = "'1.20"

# --- New row 3: Newport County vs AFC Wimbledon ---
$ws.Range("A3").Value = 'Základní Poisson'
$ws.Range("B3").Value = 'e3'
$ws.Range("C3").Value = 'Newport County'
$ws.Range("D3").Value = 'AFC Wimbledon'
$ws.Range("E3").Value = 1474.753588871106
$ws.Range("F3").Value = 1539.074720758685
$ws.Range("G3").Value = "'7.78"
$ws.Range("H3").Value = "'5.19"
$ws.Range("I3").Value = "'1.53"
$ws.Range("J3").Value = "'1.88"
$ws.Range("K3").Value = "'2.13"
$ws.Range("L3").Value = 'Bivariantní Poisson'
$ws.Range("M3").Value = "'3.59"
$ws.Range("N3").Value = "'4.00"
$ws.Range("O3").Value = "'2.12"
$ws.Range("P3").Value = "'1.70"
$ws.Range("Q3").Value = "'2.44"
$ws.Range("R3").Value = 'Monte Carlo'
$ws.Range("S3").Value = "'3.41"
$ws.Range("T3").Value = "'4.34"
$ws.Range("U3").Value = "'2.10"
$ws.Range("V3").Value = "'1.65"
$ws.Range("W3").Value = "'2.53"
$ws.Range("X3").Value = 'xGBoost'
$ws.Range("Y3").Value = "'3.99"
$ws.Range("Z3").Value = "'1.73"
$ws.Range("AA3").Value = "'5.80"
$ws.Range("AB3").Value = "'1.16"
$ws.Range("AC3").Value = "'7.42"

# The apostrophe-prefix trick above leaves a "quote prefix" style on the
# cell; clear formatting on the whole new block so it matches the plain
# (unstyled) look of the rest of the sheet.
$ws.Range("A2:AC3").ClearFormats()

# --- Existing rows (now 4-9): refresh the odds from Method 2 (Bivariantni
# Poisson), Method 3 (Monte Carlo) and Method 4 (xGBoost); Method 1 odds,
# team names and ELO numbers for these matches are unchanged. ---

# Row 4
$ws.Range("M4").Value = "'1.39"
$ws.Range("N4").Value = "'6.61"
$ws.Range("O4").Value = "'7.85"
$ws.Range("P4").Value = "'1.32"
$ws.Range("Q4").Value = "'4.14"
$ws.Range("S4").Value = "'1.40"
$ws.Range("T4").Value = "'6.60"
$ws.Range("U4").Value = "'7.34"
$ws.Range("V4").Value = "'1.30"
$ws.Range("W4").Value = "'4.30"
$ws.Range("Y4").Value = "'8.47"
$ws.Range("Z4").Value = "'1.21"
$ws.Range("AA4").Value = "'19.06"
$ws.Range("AB4").Value = "'13.88"
$ws.Range("AC4").Value = "'1.08"

# Row 5
$ws.Range("M5").Value = "'1.85"
$ws.Range("N5").Value = "'4.73"
$ws.Range("O5").Value = "'4.04"
$ws.Range("P5").Value = "'1.46"
$ws.Range("Q5").Value = "'3.19"
$ws.Range("S5").Value = "'1.87"
$ws.Range("T5").Value = "'4.75"
$ws.Range("U5").Value = "'3.91"
$ws.Range("V5").Value = "'1.44"
$ws.Range("W5").Value = "'3.27"
$ws.Range("Y5").Value = "'11.57"
$ws.Range("Z5").Value = "'1.18"
$ws.Range("AA5").Value = "'14.69"
$ws.Range("AB5").Value = "'24.28"
$ws.Range("AC5").Value = "'1.04"

# Row 6
$ws.Range("M6").Value = "'2.19"
$ws.Range("N6").Value = "'4.14"
$ws.Range("O6").Value = "'3.32"
$ws.Range("P6").Value = "'1.68"
$ws.Range("Q6").Value = "'2.47"
$ws.Range("S6").Value = "'2.20"
$ws.Range("T6").Value = "'4.20"
$ws.Range("U6").Value = "'3.26"
$ws.Range("V6").Value = "'1.64"
$ws.Range("W6").Value = "'2.55"
$ws.Range("Y6").Value = "'7.38"
$ws.Range("Z6").Value = "'1.34"
$ws.Range("AA6").Value = "'8.47"
$ws.Range("AB6").Value = "'1.33"
$ws.Range("AC6").Value = "'4.06"

# Row 7
$ws.Range("M7").Value = "'2.44"
$ws.Range("N7").Value = "'3.40"
$ws.Range("O7").Value = "'3.38"
$ws.Range("P7").Value = "'2.36"
$ws.Range("Q7").Value = "'1.74"
$ws.Range("S7").Value = "'2.38"
$ws.Range("T7").Value = "'3.67"
$ws.Range("U7").Value = "'3.25"
$ws.Range("V7").Value = "'2.33"
$ws.Range("W7").Value = "'1.75"
$ws.Range("Y7").Value = "'2.22"
$ws.Range("Z7").Value = "'6.67"
$ws.Range("AA7").Value = "'2.50"
$ws.Range("AB7").Value = "'3.26"
$ws.Range("AC7").Value = "'1.44"

# Row 8
$ws.Range("M8").Value = "'2.33"
$ws.Range("N8").Value = "'4.61"
$ws.Range("O8").Value = "'2.83"
$ws.Range("P8").Value = "'1.36"
$ws.Range("Q8").Value = "'3.79"
$ws.Range("S8").Value = "'2.25"
$ws.Range("T8").Value = "'4.96"
$ws.Range("U8").Value = "'2.82"
$ws.Range("V8").Value = "'1.34"
$ws.Range("W8").Value = "'3.93"
$ws.Range("Y8").Value = "'7.32"
$ws.Range("Z8").Value = "'6.00"
$ws.Range("AA8").Value = "'1.44"
$ws.Range("AB8").Value = "'1.07"
$ws.Range("AC8").Value = "'15.42"

# Row 9
$ws.Range("M9").Value = "'2.38"
$ws.Range("N9").Value = "'4.79"
$ws.Range("O9").Value = "'2.70"
$ws.Range("P9").Value = "'1.29"
$ws.Range("Q9").Value = "'4.42"
$ws.Range("S9").Value = "'2.39"
$ws.Range("T9").Value = "'4.97"
$ws.Range("U9").Value = "'2.63"
$ws.Range("V9").Value = "'1.26"
$ws.Range("W9").Value = "'4.85"
$ws.Range("Y9").Value = "'7.10"
$ws.Range("Z9").Value = "'1.22"
$ws.Range("AA9").Value = "'23.64"
$ws.Range("AB9").Value = "'1.29"
$ws.Range("AC9").Value = "'4.44"

# Again, strip the quote-prefix style picked up from the apostrophe trick
# on each of the three updated odds blocks (Method 2 / 3 / 4 columns).
$ws.Range("M4:Q9").ClearFormats()
$ws.Range("S4:W9").ClearFormats()
$ws.Range("Y4:AC9").ClearFormats()
